# Separate Sidebars as per roles
# Adds a new "Pages" worksheet (after the existing "Roles" sheet), makes it
# the active/selected tab, and records a page route in cell A1.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet ("Roles") so it lands
# at the end of the tab strip, then rename it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Pages"

# Content
$ws.Range("A1").Value = "ecommerce-page/order-details"

# Column sizing / selection to match the authored sheet
$ws.Columns("A").ColumnWidth = 39.94
$ws.Range("A8").Select() | Out-Null
